$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from right after "MP73010" (in the
#    title line) to the very end of the third paragraph, right after
#    "... assignment report!" (and before that paragraph's mark).
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# Find the end of paragraph 3 (the paragraph ending in "assignment report!"),
# i.e. the position right before that paragraph's own paragraph mark.
$para3 = $d.Paragraphs(3)
$insertionPoint = $para3.Range.End - 1

# Directly adding a bookmark exactly one character before a paragraph
# mark can misplace it, so nudge the text out of the way first: insert
# a one-character placeholder, drop the bookmark in front of it (a
# location that is NOT directly adjacent to the paragraph mark), then
# remove the placeholder again.
$placeholder = $d.Range($insertionPoint, $insertionPoint)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($insertionPoint, $insertionPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($insertionPoint, $insertionPoint + 1).Delete()

# ------------------------------------------------------------------
# 2) Add a new paragraph after "Ben changing things up!" containing:
#    "Instructions Noted Down ,Let us do this !!" with gramStart/
#    gramEnd proofErr markers bracketing "Down ,Let".
# ------------------------------------------------------------------
$benPara = $d.Paragraphs(5)
$afterBen = $benPara.Range.Duplicate
$afterBen.Collapse(0)
$afterBen.InsertParagraphAfter()

$newPara = $d.Paragraphs(6)
$target = $d.Range($newPara.Range.Start, $newPara.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Instructions Noted </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Down ,Let</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> us do this !!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
